# Apply "Moderhinke Dashboard.xlsx" weekly refresh:
# 1) Revise several already-reported farms_total_count /
#    farms_to_examine_count figures downward (late-arriving corrections
#    for the most recent weeks already present in the sheet).
# 2) Append the new ISO week 202515 (week ending 2025-04-13) block of
#    5 rows (farms_total_count, farms_to_examine_count, farms_examined_count,
#    farms_examined_positive_count, farms_examined_negative_count).
# 3) Update the sheet's scroll/selection state to match where the author
#    left the cursor after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Revised counts for already-present weeks -------------------------
$ws.Range("D2").Value = 11721
$ws.Range("D3").Value = 11466
$ws.Range("D7").Value = 11826
$ws.Range("D8").Value = 11113
$ws.Range("D12").Value = 11934
$ws.Range("D13").Value = 10611
$ws.Range("D17").Value = 11978
$ws.Range("D18").Value = 10060
$ws.Range("D22").Value = 12036
$ws.Range("D23").Value = 9632
$ws.Range("D27").Value = 12077
$ws.Range("D28").Value = 9172
$ws.Range("D32").Value = 12117
$ws.Range("D33").Value = 8764
$ws.Range("D37").Value = 12155
$ws.Range("D38").Value = 8240
$ws.Range("D42").Value = 12190
$ws.Range("D43").Value = 7673
$ws.Range("D47").Value = 12215
$ws.Range("D48").Value = 7125
$ws.Range("D52").Value = 12243
$ws.Range("D53").Value = 6483
$ws.Range("D57").Value = 12266
$ws.Range("D58").Value = 5872
$ws.Range("D62").Value = 12277
$ws.Range("D63").Value = 5706
$ws.Range("D67").Value = 12299
$ws.Range("D68").Value = 5522
$ws.Range("D72").Value = 12315
$ws.Range("D73").Value = 5086
$ws.Range("D77").Value = 12337
$ws.Range("D78").Value = 4608
$ws.Range("D82").Value = 12355
$ws.Range("D83").Value = 4089
$ws.Range("D87").Value = 12382
$ws.Range("D88").Value = 3683
$ws.Range("D92").Value = 12397
$ws.Range("D93").Value = 3344
$ws.Range("D97").Value = 12418
$ws.Range("D98").Value = 3033
$ws.Range("D102").Value = 12433
$ws.Range("D103").Value = 2757
$ws.Range("D107").Value = 12443
$ws.Range("D108").Value = 2470
$ws.Range("D112").Value = 12464
$ws.Range("D113").Value = 2233
$ws.Range("D117").Value = 12485
$ws.Range("D118").Value = 1988
$ws.Range("D122").Value = 12502
$ws.Range("D123").Value = 1734
$ws.Range("D127").Value = 12536
$ws.Range("D128").Value = 1475
$ws.Range("D132").Value = 12553
$ws.Range("D133").Value = 1283

# --- 2) Append new week 202515 (LastDayOfWeek 2025-04-13 / serial 45760) -
$ws.Range("A137").Value = 202515
$ws.Range("B137").Value = 45760
$ws.Range("C137").Value = "farms_total_count"
$ws.Range("D137").Value = 12565

$ws.Range("A138").Value = 202515
$ws.Range("B138").Value = 45760
$ws.Range("C138").Value = "farms_to_examine_count"
$ws.Range("D138").Value = 1255

$ws.Range("A139").Value = 202515
$ws.Range("B139").Value = 45760
$ws.Range("C139").Value = "farms_examined_count"
$ws.Range("D139").Value = 11310

$ws.Range("A140").Value = 202515
$ws.Range("B140").Value = 45760
$ws.Range("C140").Value = "farms_examined_positive_count"
$ws.Range("D140").Value = 1381

$ws.Range("A141").Value = 202515
$ws.Range("B141").Value = 45760
$ws.Range("C141").Value = "farms_examined_negative_count"
$ws.Range("D141").Value = 9929

# The B column uses a date display format (m/d/yyyy, same style as the
# rest of column B) - copy that formatting from the row above onto each
# of the freshly written B cells instead of re-declaring a numeric format.
$ws.Range("B136").Copy()
$ws.Range("B137:B141").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 3) Leave the view where the author left it after entering the data --
$ws.Range("G134").Select()
